$wb = $excel.ActiveWorkbook

# Each localized sheet stamps column G ("Latest Handback DateTime") with
# the report-generation timestamp for every row that has been handed
# back. Regenerating the report bumps that timestamp forward by about a
# minute on each of the four localization sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G1:G271").Replace("2016-02-22 08:49:28", "2016-02-22 08:50:29")

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("G1:G271").Replace("2016-02-22 08:49:38", "2016-02-22 08:50:40")

$wsJaJp = $wb.Worksheets.Item("ja-jp")
$wsJaJp.Range("G1:G271").Replace("2016-02-22 08:49:48", "2016-02-22 08:50:50")

$wsZhTw = $wb.Worksheets.Item("zh-tw")
$wsZhTw.Range("G1:G271").Replace("2016-02-22 08:49:58", "2016-02-22 08:51:01")
